# Auto-generated: apply scheduled market-price refresh to Carbuncle_Profits workbook
# Updates cached currentAveragePrice / LevePrice / LeveProfit columns (H-N) per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29: H29
$ws.Range("H29").Value = 1483.3334
# Row 64: H64,I64,J64,K64,L64,M64,N64
$ws.Range("H64").Value = 4175.6343
$ws.Range("I64").Value = 3318.75
$ws.Range("J64").Value = 7222.3335
$ws.Range("K64").Value = 3318.75
$ws.Range("L64").Value = 7222.3335
$ws.Range("M64").Value = -3070.75
$ws.Range("N64").Value = -7718.3335
# Row 67: H67,I67,J67,K67,L67,M67,N67
$ws.Range("H67").Value = 4175.6343
$ws.Range("I67").Value = 3318.75
$ws.Range("J67").Value = 7222.3335
$ws.Range("K67").Value = 3318.75
$ws.Range("L67").Value = 7222.3335
$ws.Range("M67").Value = -2460.75
$ws.Range("N67").Value = -8938.333500000001
# Row 87: H87,J87,L87,N87
$ws.Range("H87").Value = 19448.906
$ws.Range("J87").Value = 19448.906
$ws.Range("L87").Value = 19448.906
$ws.Range("N87").Value = -21944.906
# Row 90: H90,J90,L90,N90
$ws.Range("H90").Value = 19448.906
$ws.Range("J90").Value = 19448.906
$ws.Range("L90").Value = 58346.71799999999
$ws.Range("N90").Value = -70826.71799999999
# Row 129: H129,I129,J129,K129,L129,M129,N129
$ws.Range("H129").Value = 1061.7354
$ws.Range("I129").Value = 412
$ws.Range("J129").Value = 1295.64
$ws.Range("K129").Value = 1236
$ws.Range("L129").Value = 3886.92
$ws.Range("M129").Value = 3764
$ws.Range("N129").Value = -13886.92
# Row 135: H135,I135,K135,M135
$ws.Range("H135").Value = 460.96875
$ws.Range("I135").Value = 425.89655
$ws.Range("K135").Value = 3833.06895
$ws.Range("M135").Value = -1298.06895
# Row 137: H137,I137,J137,K137,L137,M137,N137
$ws.Range("H137").Value = 314298.16
$ws.Range("I137").Value = 477311.66
$ws.Range("J137").Value = 3090.5454
$ws.Range("K137").Value = 1431934.98
$ws.Range("L137").Value = 9271.636200000001
$ws.Range("M137").Value = -1429384.98
$ws.Range("N137").Value = -14371.6362
# Row 138: H138,I138,J138,K138,L138,M138,N138
$ws.Range("H138").Value = 2861.4546
$ws.Range("I138").Value = 1314.3077
$ws.Range("J138").Value = 6632.625
$ws.Range("K138").Value = 3942.9231
$ws.Range("L138").Value = 19897.875
$ws.Range("M138").Value = 1197.0769
$ws.Range("N138").Value = -30177.875
# Row 141: H141,I141,K141,M141
$ws.Range("H141").Value = 9906.781000000001
$ws.Range("I141").Value = 1471.6471
$ws.Range("K141").Value = 4414.9413
$ws.Range("M141").Value = 765.0587000000005

$ws = $wb.Worksheets.Item("ARM")
# Row 31: H31,I31,K31,M31
$ws.Range("H31").Value = 8985.5
$ws.Range("I31").Value = 8985.5
$ws.Range("K31").Value = 8985.5
$ws.Range("M31").Value = -8691.5
# Row 32: H32,I32,J32,K32,L32,M32,N32
$ws.Range("H32").Value = 4938.6216
$ws.Range("I32").Value = 3673.6365
$ws.Range("J32").Value = 15374.75
$ws.Range("K32").Value = 3673.6365
$ws.Range("L32").Value = 15374.75
$ws.Range("M32").Value = -3386.6365
$ws.Range("N32").Value = -15948.75
# Row 61: H61,I61,J61,K61,L61,M61,N61
$ws.Range("H61").Value = 1508.3889
$ws.Range("I61").Value = 902.7143
$ws.Range("J61").Value = 3628.25
$ws.Range("K61").Value = 902.7143
$ws.Range("L61").Value = 3628.25
$ws.Range("M61").Value = -690.7143
$ws.Range("N61").Value = -4052.25
# Row 88: H88,I88,J88,K88,L88,M88,N88
$ws.Range("H88").Value = 2000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2000
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -2812
# Row 91: H91,I91,J91,K91,L91,M91,N91
$ws.Range("H91").Value = 2000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2000
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -4808
# Row 110: H110,I110,K110,M110
$ws.Range("H110").Value = 74229.71000000001
$ws.Range("I110").Value = 102125.6
$ws.Range("K110").Value = 102125.6
$ws.Range("M110").Value = -100080.6
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 1508.3889
$ws.Range("I136").Value = 902.7143
$ws.Range("J136").Value = 3628.25
$ws.Range("K136").Value = 2708.1429
$ws.Range("L136").Value = 10884.75
$ws.Range("M136").Value = -158.1428999999998
$ws.Range("N136").Value = -15984.75

$ws = $wb.Worksheets.Item("BSM")
# Row 86: H86,I86,J86,K86,L86,M86,N86
$ws.Range("H86").Value = 3003.8333
$ws.Range("I86").Value = 3258.25
$ws.Range("J86").Value = 2495
$ws.Range("K86").Value = 3258.25
$ws.Range("L86").Value = 2495
$ws.Range("M86").Value = -2135.25
$ws.Range("N86").Value = -4741
# Row 89: H89,I89,J89,K89,L89,M89,N89
$ws.Range("H89").Value = 3003.8333
$ws.Range("I89").Value = 3258.25
$ws.Range("J89").Value = 2495
$ws.Range("K89").Value = 16291.25
$ws.Range("L89").Value = 12475
$ws.Range("M89").Value = -10675.25
$ws.Range("N89").Value = -23707
# Row 107: H107,I107,J107,K107,L107,M107,N107
$ws.Range("H107").Value = 806.9524
$ws.Range("I107").Value = 786.2778
$ws.Range("J107").Value = 931
$ws.Range("K107").Value = 786.2778
$ws.Range("L107").Value = 931
$ws.Range("M107").Value = 1133.7222
$ws.Range("N107").Value = -4771

$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31,I31,J31,K31,L31,M31,N31
$ws.Range("H31").Value = 2589.0176
$ws.Range("I31").Value = 2391.5
$ws.Range("J31").Value = 2680.1794
$ws.Range("K31").Value = 2391.5
$ws.Range("L31").Value = 2680.1794
$ws.Range("M31").Value = -2096.5
$ws.Range("N31").Value = -3270.1794
# Row 34: H34,I34,J34,K34,L34,M34,N34
$ws.Range("H34").Value = 2589.0176
$ws.Range("I34").Value = 2391.5
$ws.Range("J34").Value = 2680.1794
$ws.Range("K34").Value = 2391.5
$ws.Range("L34").Value = 2680.1794
$ws.Range("M34").Value = -2189.5
$ws.Range("N34").Value = -3084.1794
# Row 62: H62,I62,J62,K62,L62,M62,N62
$ws.Range("H62").Value = 4709.4546
$ws.Range("I62").Value = 2983.1667
$ws.Range("J62").Value = 6781
$ws.Range("K62").Value = 2983.1667
$ws.Range("L62").Value = 6781
$ws.Range("M62").Value = -2359.1667
$ws.Range("N62").Value = -8029
# Row 65: H65,I65,J65,K65,L65,M65,N65
$ws.Range("H65").Value = 4709.4546
$ws.Range("I65").Value = 2983.1667
$ws.Range("J65").Value = 6781
$ws.Range("K65").Value = 14915.8335
$ws.Range("L65").Value = 33905
$ws.Range("M65").Value = -11795.8335
$ws.Range("N65").Value = -40145
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 2341.6667
$ws.Range("I132").Value = 967.65216
$ws.Range("J132").Value = 4772.615
$ws.Range("K132").Value = 2902.95648
$ws.Range("L132").Value = 14317.845
$ws.Range("M132").Value = -372.9564799999998
$ws.Range("N132").Value = -19377.845
# Row 134: H134,I134,J134,K134,L134,M134,N134
$ws.Range("H134").Value = 2756.5908
$ws.Range("I134").Value = 2431.55
$ws.Range("J134").Value = 6007
$ws.Range("K134").Value = 7294.650000000001
$ws.Range("L134").Value = 18021
$ws.Range("M134").Value = -4759.650000000001
$ws.Range("N134").Value = -23091

$ws = $wb.Worksheets.Item("CUL")
# Row 131: H131,J131,L131,N131
$ws.Range("H131").Value = 4479.6206
$ws.Range("J131").Value = 7788.0625
$ws.Range("L131").Value = 23364.1875
$ws.Range("N131").Value = -33444.1875
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 1999.75
$ws.Range("I136").Value = 2054.9167
$ws.Range("J136").Value = 1979.0625
$ws.Range("K136").Value = 6164.750100000001
$ws.Range("L136").Value = 5937.1875
$ws.Range("M136").Value = -1064.750100000001
$ws.Range("N136").Value = -16137.1875

$ws = $wb.Worksheets.Item("LTW")
# Row 22: H22,J22,L22,N22
$ws.Range("H22").Value = 1365.3334
$ws.Range("J22").Value = 1764.6666
$ws.Range("L22").Value = 1764.6666
$ws.Range("N22").Value = -2354.6666
# Row 27: H27,J27,L27,N27
$ws.Range("H27").Value = 1365.3334
$ws.Range("J27").Value = 1764.6666
$ws.Range("L27").Value = 1764.6666
$ws.Range("N27").Value = -1978.6666
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 13335071
$ws.Range("I136").Value = 1809.6666
$ws.Range("J136").Value = 333333340
$ws.Range("K136").Value = 5428.9998
$ws.Range("L136").Value = 1000000020
$ws.Range("M136").Value = -2878.9998
$ws.Range("N136").Value = -1000005120

$ws = $wb.Worksheets.Item("WVR")
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 5378199.5
$ws.Range("I136").Value = 7247353.5
$ws.Range("J136").Value = 4381.5625
$ws.Range("K136").Value = 21742060.5
$ws.Range("L136").Value = 13144.6875
$ws.Range("M136").Value = -21739510.5
$ws.Range("N136").Value = -18244.6875
